$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = '27.504.05'
$d.ClearFormats()
$ws.Range("E2").Value = '  +2.22%  '

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = '1.841.70'
$d.ClearFormats()
$ws.Range("E3").Value = '  +1.56%  '

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = '1.016'
$d.ClearFormats()
$ws.Range("E4").Value = '  +1.41%  '

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = '315.37'
$d.ClearFormats()
$ws.Range("E5").Value = '  +2.09%  '

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = '1.013'
$d.ClearFormats()
$ws.Range("E6").Value = '  +1.15%  '

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = '0.4757'
$d.ClearFormats()
$ws.Range("E7").Value = '  +2.20%  '

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = '0.3706'
$d.ClearFormats()
$ws.Range("E8").Value = '  +1.07%  '

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = '0.07474'
$d.ClearFormats()
$ws.Range("E9").Value = '  +1.66%  '

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = '0.8866'
$d.ClearFormats()
$ws.Range("E10").Value = '  +2.23%  '

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = '20.53'
$d.ClearFormats()
$ws.Range("E11").Value = '  +1.21%  '

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = '1.876.59'
$d.ClearFormats()
$ws.Range("E12").Value = '  +3.36%  '

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = '0.07373'
$d.ClearFormats()
$ws.Range("E13").Value = '  +4.12%  '

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = '5.468'
$d.ClearFormats()
$ws.Range("E14").Value = '  +1.77%  '

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = '93.25'
$d.ClearFormats()
$ws.Range("E15").Value = '  +1.92%  '

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = '6.594'
$d.ClearFormats()
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("E17").Value = '  +1.12%  '

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = '0.000008836'
$d.ClearFormats()
$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("E19").Value = '  +1.17%  '

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = '14.83'
$d.ClearFormats()
$ws.Range("E20").Value = '  +1.20%  '

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = '27.527.81'
$d.ClearFormats()
$ws.Range("E21").Value = '  +2.22%  '

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = '5.331'
$d.ClearFormats()
$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("E23").Value = '  +1.07%  '

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = '2.105.23'
$d.ClearFormats()
$ws.Range("E24").Value = '  +3.07%  '

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = '1.900'
$d.ClearFormats()
$ws.Range("E25").Value = '  +0.40%  '

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = '152.24'
$d.ClearFormats()
$ws.Range("E26").Value = '  +1.39%  '

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = '18.64'
$d.ClearFormats()
$ws.Range("E27").Value = '  +2.00%  '

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = '2.157'
$d.ClearFormats()
$ws.Range("E28").Value = '  -0.36%  '

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = '5.270'
$d.ClearFormats()
$ws.Range("E29").Value = '  +0.10%  '

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = '118.27'
$d.ClearFormats()
$ws.Range("E30").Value = '  +2.48%  '

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = '0.09006'
$d.ClearFormats()
$ws.Range("E31").Value = '  +0.75%  '

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = '0.7575'
$d.ClearFormats()
$ws.Range("E32").Value = '  +0.02%  '

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = '1.181'
$d.ClearFormats()
$ws.Range("E33").Value = '  +2.06%  '

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = '4.565'
$d.ClearFormats()
$ws.Range("E34").Value = '  +1.73%  '

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = '2.970'
$d.ClearFormats()
$ws.Range("E35").Value = '  +2.00%  '

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = '1.014'
$d.ClearFormats()
$ws.Range("E36").Value = '  +1.28%  '

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = '1.108'
$d.ClearFormats()
$ws.Range("E37").Value = '  +2.30%  '

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = '0.05348'
$d.ClearFormats()
$ws.Range("E38").Value = '  +1.24%  '

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = '0.01962'
$d.ClearFormats()
$ws.Range("E39").Value = '  +0.52%  '

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = '3.007'
$d.ClearFormats()
$ws.Range("E40").Value = '  +0.16%  '

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = '7.330'
$d.ClearFormats()
$ws.Range("E41").Value = '  +1.91%  '

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = '2.392'
$d.ClearFormats()
$ws.Range("E42").Value = '  +4.43%  '

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = '0.5350'
$d.ClearFormats()
$ws.Range("E43").Value = '  +0.96%  '

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = '0.1664'
$d.ClearFormats()
$ws.Range("E44").Value = '  +0.75%  '

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = '8.545'
$d.ClearFormats()
$ws.Range("E45").Value = '  +1.88%  '

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = '0.4946'
$d.ClearFormats()
$ws.Range("E46").Value = '  +1.74%  '

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = '10.59'
$d.ClearFormats()
$ws.Range("E47").Value = '  +1.68%  '

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = '1.015'
$d.ClearFormats()
$ws.Range("E48").Value = '  +1.33%  '

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = '104.84'
$d.ClearFormats()
$ws.Range("E49").Value = '  +2.02%  '

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = '1.683'
$d.ClearFormats()
$ws.Range("E50").Value = '  +1.44%  '

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = '0.06327'
$d.ClearFormats()
$ws.Range("E51").Value = '  +0.65%  '

